$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain-text Price value into a cell without Excel silently
# re-interpreting it as a number (which would drop formatting like trailing
# zeros). We temporarily force Text format, assign the value, then restore the
# cell to the default (unstyled) "Normal" style so no stray formatting remains.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "61.417.73"
$ws.Range("E2").Value = "  +8.40%  "

# Row 3
$ws.Range("D3").Value = "2.671.18"
$ws.Range("E3").Value = "  +10.42%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
Set-TextValue $ws.Range("D5") "515.79"
$ws.Range("E5").Value = "  +5.87%  "

# Row 6
Set-TextValue $ws.Range("D6") "160.40"
$ws.Range("E6").Value = "  +5.10%  "

# Row 7
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.43%  "

# Row 8
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws.Range("D8") "0.615"
$ws.Range("E8").Value = "  +1.33%  "

# Row 9
$ws.Range("D9").Value = "2.669.60"
$ws.Range("E9").Value = "  +9.56%  "

# Row 10
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D10") "6.14"
$ws.Range("E10").Value = "  +10.00%  "

# Row 11
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D11") "0.106"
$ws.Range("E11").Value = "  +7.34%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.351"
$ws.Range("E12").Value = "  +5.10%  "

# Row 13
$ws.Range("E13").Value = "  +1.57%  "

# Row 14
$ws.Range("D14").Value = "3.131.75"
$ws.Range("E14").Value = "  +10.52%  "

# Row 15
$ws.Range("D15").Value = "61.210.54"
$ws.Range("E15").Value = "  +8.72%  "

# Row 16
Set-TextValue $ws.Range("D16") "22.44"
$ws.Range("E16").Value = "  +9.47%  "

# Row 17
$ws.Range("E17").Value = "  +6.67%  "

# Row 18
$ws.Range("D18").Value = "2.666.00"
$ws.Range("E18").Value = "  +10.06%  "

# Row 19
Set-TextValue $ws.Range("D19") "4.86"
$ws.Range("E19").Value = "  +3.66%  "

# Row 20
Set-TextValue $ws.Range("D20") "354.02"
$ws.Range("E20").Value = "  +9.27%  "

# Row 21
Set-TextValue $ws.Range("D21") "10.59"
$ws.Range("E21").Value = "  +7.76%  "

# Row 22
Set-TextValue $ws.Range("D22") "6.22"
$ws.Range("E22").Value = "  +6.06%  "

# Row 23
Set-TextValue $ws.Range("D23") "1.00"
$ws.Range("E23").Value = "  +0.13%  "

# Row 24
Set-TextValue $ws.Range("D24") "60.80"
$ws.Range("E24").Value = "  +5.94%  "

# Row 25
Set-TextValue $ws.Range("D25") "0.429"
$ws.Range("E25").Value = "  +5.99%  "

# Row 26
$ws.Range("D26").Value = "2.781.87"
$ws.Range("E26").Value = "  +11.03%  "

# Row 27
Set-TextValue $ws.Range("D27") "0.168"
$ws.Range("E27").Value = "  +6.04%  "

# Row 28
Set-TextValue $ws.Range("D28") "0.999"
$ws.Range("E28").Value = "  +0.48%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0878"
$ws.Range("E29").Value = "  +14.51%  "

# Row 30
Set-TextValue $ws.Range("D30") "7.55"
$ws.Range("E30").Value = "  +2.69%  "

# Row 31
$ws.Range("E31").Value = "  +0.23%  "

# Row 32
Set-TextValue $ws.Range("D32") "19.78"
$ws.Range("E32").Value = "  +7.06%  "

# Row 33
Set-TextValue $ws.Range("D33") "157.36"
$ws.Range("E33").Value = "  +5.69%  "

# Row 34
$ws.Range("E34").Value = "  +5.62%  "

# Row 35
Set-TextValue $ws.Range("D35") "5.78"
$ws.Range("E35").Value = "  +10.47%  "

# Row 36
Set-TextValue $ws.Range("D36") "4.13"
$ws.Range("E36").Value = "  +12.43%  "

# Row 37
$ws.Range("E37").Value = "  +9.06%  "

# Row 38
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D38") "1.56"
$ws.Range("E38").Value = "  +14.03%  "

# Row 39
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D39") "0.884"
$ws.Range("E39").Value = "  +4.81%  "

# Row 40
Set-TextValue $ws.Range("D40") "3.80"
$ws.Range("E40").Value = "  +9.45%  "

# Row 41
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D41") "303.47"
$ws.Range("E41").Value = "  +16.09%  "

# Row 42
$ws.Range("B42").Value = "SuiNetwork"
$ws.Range("C42").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D42") "0.838"
$ws.Range("E42").Value = "  +34.44%  "

# Row 43
Set-TextValue $ws.Range("D43") "35.80"
$ws.Range("E43").Value = "  +4.57%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.648"
$ws.Range("E44").Value = "  +8.17%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.0580"
$ws.Range("E45").Value = "  +9.31%  "

# Row 46
$ws.Range("E46").Value = "  +1.50%  "

# Row 47
Set-TextValue $ws.Range("D47") "20.28"
$ws.Range("E47").Value = "  +17.07%  "

# Row 48
Set-TextValue $ws.Range("D48") "0.999"
$ws.Range("E48").Value = "  +0.48%  "

# Row 49
Set-TextValue $ws.Range("D49") "5.01"
$ws.Range("E49").Value = "  +6.69%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.0239"
$ws.Range("E50").Value = "  +5.29%  "

# Row 51
$ws.Range("D51").Value = "2.031.86"
$ws.Range("E51").Value = "  +9.74%  "
